$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for column C, matching style of existing headers (A1/B1)
$ws.Range("C1").Value = "total_points"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Fill in points values for rows 2-5 (0 for first three customers, 30 for the updated one)
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 30

# Cell A5 was stored as text "79172233" - normalize it to a real number like the diff expects
$ws.Range("A5").Value = 79172233
